$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.214.50"
$ws.Range("E2").Value = "  +2.62%  "
$ws.Range("D3").Value = "3.600.95"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'608.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.61%  "
$ws.Range("D6").Value = "'175.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.79%  "
$ws.Range("D7").Value = "3.594.76"
$ws.Range("E7").Value = "  +1.72%  "
$ws.Range("D8").Value = "'0.620"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.34%  "
$ws.Range("D10").Value = "'0.199"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.80%  "
$ws.Range("D11").Value = "'7.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +11.04%  "
$ws.Range("D12").Value = "'0.593"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("D13").Value = "'47.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("D15").Value = "4.167.91"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").Value = "'8.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("D17").Value = "'625.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").Value = "3.585.11"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("D19").Value = "71.246.03"
$ws.Range("E19").Value = "  +2.54%  "
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("D21").Value = "'17.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").Value = "'0.894"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("E23").Value = "  -15.89%  "
$ws.Range("D24").Value = "'16.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.07%  "
$ws.Range("D25").Value = "'98.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").Value = "'3.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "'2.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("D29").Value = "'9.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("D30").Value = "'33.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.93%  "
$ws.Range("D31").Value = "'8.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").Value = "'3.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("E34").Value = "  -2.59%  "
$ws.Range("D35").Value = "'644.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("D36").Value = "'3.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.24%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'10.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("D39").Value = "'0.0487"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.69%  "
$ws.Range("D40").Value = "'57.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("E42").Value = "  +4.30%  "
$ws.Range("D43").Value = "3.418.78"
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("D44").Value = "'0.327"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("D45").Value = "0.0₃0724"
$ws.Range("E45").Value = "  +2.99%  "
$ws.Range("E46").Value = "  +9.11%  "
$ws.Range("D47").Value = "'33.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.53%  "
$ws.Range("D48").Value = "'2.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.17%  "
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("D50").Value = "'133.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("E51").Value = "  -0.02%  "
